$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 54.2
$ws.Range("J9").Value = 63.333332
$ws.Range("L9").Value = 63.333332
$ws.Range("N9").Value = -401.333332
$ws.Range("H20").Value = 50010.5
$ws.Range("I20").Value = 50010.5
$ws.Range("K20").Value = 50010.5
$ws.Range("M20").Value = -49780.5
$ws.Range("H35").Value = 50010.5
$ws.Range("I35").Value = 50010.5
$ws.Range("K35").Value = 50010.5
$ws.Range("M35").Value = -49631.5
$ws.Range("H129").Value = 1423.4736
$ws.Range("J129").Value = 1634.4193
$ws.Range("L129").Value = 4903.257900000001
$ws.Range("N129").Value = -14903.2579
$ws.Range("H132").Value = 445675.6
$ws.Range("I132").Value = 507437.72
$ws.Range("J132").Value = 75103
$ws.Range("K132").Value = 1522313.16
$ws.Range("L132").Value = 225309
$ws.Range("M132").Value = -1519783.16
$ws.Range("N132").Value = -230369
$ws.Range("H137").Value = 22223428
$ws.Range("I137").Value = 38462492
$ws.Range("J137").Value = 1551.4736
$ws.Range("K137").Value = 115387476
$ws.Range("L137").Value = 4654.4208
$ws.Range("M137").Value = -115384926
$ws.Range("N137").Value = -9754.4208
$ws.Range("H138").Value = 1772.86
$ws.Range("I138").Value = 586.5282999999999
$ws.Range("J138").Value = 3110.6382
$ws.Range("K138").Value = 1759.5849
$ws.Range("L138").Value = 9331.9146
$ws.Range("M138").Value = 3380.4151
$ws.Range("N138").Value = -19611.9146
$ws.Range("H141").Value = 1921.9143
$ws.Range("I141").Value = 1221.85
$ws.Range("J141").Value = 6122.3
$ws.Range("K141").Value = 3665.55
$ws.Range("L141").Value = 18366.9
$ws.Range("M141").Value = 1514.45
$ws.Range("N141").Value = -28726.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H32").Value = 16947.734
$ws.Range("I32").Value = 2145.123
$ws.Range("J32").Value = 337671
$ws.Range("K32").Value = 2145.123
$ws.Range("L32").Value = 337671
$ws.Range("M32").Value = -1858.123
$ws.Range("N32").Value = -338245
$ws.Range("H61").Value = 2670.775
$ws.Range("I61").Value = 1646.5927
$ws.Range("J61").Value = 4797.923
$ws.Range("K61").Value = 1646.5927
$ws.Range("L61").Value = 4797.923
$ws.Range("M61").Value = -1434.5927
$ws.Range("N61").Value = -5221.923
$ws.Range("H74").Value = 11203.429
$ws.Range("I74").Value = 2737.3333
$ws.Range("J74").Value = 62000
$ws.Range("K74").Value = 2737.3333
$ws.Range("L74").Value = 62000
$ws.Range("M74").Value = -1863.3333
$ws.Range("N74").Value = -63748
$ws.Range("H77").Value = 11203.429
$ws.Range("I77").Value = 2737.3333
$ws.Range("J77").Value = 62000
$ws.Range("K77").Value = 13686.6665
$ws.Range("L77").Value = 310000
$ws.Range("M77").Value = -9318.666499999999
$ws.Range("N77").Value = -318736
$ws.Range("H122").Value = 2216.2
$ws.Range("I122").Value = 2066.238
$ws.Range("K122").Value = 6198.714
$ws.Range("M122").Value = -3748.714
$ws.Range("H132").Value = 2098.6216
$ws.Range("I132").Value = 1709
$ws.Range("J132").Value = 5313
$ws.Range("K132").Value = 5127
$ws.Range("L132").Value = 15939
$ws.Range("M132").Value = -2597
$ws.Range("N132").Value = -20999
$ws.Range("H136").Value = 2670.775
$ws.Range("I136").Value = 1646.5927
$ws.Range("J136").Value = 4797.923
$ws.Range("K136").Value = 4939.7781
$ws.Range("L136").Value = 14393.769
$ws.Range("M136").Value = -2389.7781
$ws.Range("N136").Value = -19493.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3376.1538
$ws.Range("I20").Value = 3387.7778
$ws.Range("J20").Value = 3350
$ws.Range("K20").Value = 3387.7778
$ws.Range("L20").Value = 3350
$ws.Range("M20").Value = -3140.7778
$ws.Range("N20").Value = -3844
$ws.Range("H86").Value = 7040.65
$ws.Range("I86").Value = 2613.25
$ws.Range("J86").Value = 9992.25
$ws.Range("K86").Value = 2613.25
$ws.Range("L86").Value = 9992.25
$ws.Range("M86").Value = -1490.25
$ws.Range("N86").Value = -12238.25
$ws.Range("H89").Value = 7040.65
$ws.Range("I89").Value = 2613.25
$ws.Range("J89").Value = 9992.25
$ws.Range("K89").Value = 13066.25
$ws.Range("L89").Value = 49961.25
$ws.Range("M89").Value = -7450.25
$ws.Range("N89").Value = -61193.25
$ws.Range("H107").Value = 884.0833
$ws.Range("I107").Value = 810.9
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 810.9
$ws.Range("L107").Value = 1250
$ws.Range("M107").Value = 1109.1
$ws.Range("N107").Value = -5090
$ws.Range("H134").Value = 2524.3274
$ws.Range("I134").Value = 1593.2439
$ws.Range("K134").Value = 4779.7317
$ws.Range("M134").Value = -2244.7317

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1882.9459
$ws.Range("I31").Value = 1083.4
$ws.Range("J31").Value = 2428.0908
$ws.Range("K31").Value = 1083.4
$ws.Range("L31").Value = 2428.0908
$ws.Range("M31").Value = -788.4000000000001
$ws.Range("N31").Value = -3018.0908
$ws.Range("H34").Value = 1882.9459
$ws.Range("I34").Value = 1083.4
$ws.Range("J34").Value = 2428.0908
$ws.Range("K34").Value = 1083.4
$ws.Range("L34").Value = 2428.0908
$ws.Range("M34").Value = -881.4000000000001
$ws.Range("N34").Value = -2832.0908
$ws.Range("H122").Value = 2296.0527
$ws.Range("I122").Value = 1242.4
$ws.Range("J122").Value = 3466.7778
$ws.Range("K122").Value = 3727.2
$ws.Range("L122").Value = 10400.3334
$ws.Range("M122").Value = -1277.2
$ws.Range("N122").Value = -15300.3334
$ws.Range("H132").Value = 1559.6578
$ws.Range("I132").Value = 1154.6
$ws.Range("J132").Value = 3953.182
$ws.Range("K132").Value = 3463.8
$ws.Range("L132").Value = 11859.546
$ws.Range("M132").Value = -933.7999999999997
$ws.Range("N132").Value = -16919.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1207.8889
$ws.Range("J107").Value = 601.7143
$ws.Range("L107").Value = 1805.1429
$ws.Range("N107").Value = -5645.1429
$ws.Range("H131").Value = 1722.579
$ws.Range("I131").Value = 382.85715
$ws.Range("J131").Value = 2025.0968
$ws.Range("K131").Value = 1148.57145
$ws.Range("L131").Value = 6075.2904
$ws.Range("M131").Value = 3891.42855
$ws.Range("N131").Value = -16155.2904
$ws.Range("H133").Value = 4907.615
$ws.Range("I133").Value = 2779.9
$ws.Range("K133").Value = 8339.700000000001
$ws.Range("M133").Value = -3279.700000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 856282.6
$ws.Range("I122").Value = 2223577
$ws.Range("J122").Value = 1723.625
$ws.Range("K122").Value = 6670731
$ws.Range("L122").Value = 5170.875
$ws.Range("M122").Value = -6668281
$ws.Range("N122").Value = -10070.875
$ws.Range("H132").Value = 2951.4634
$ws.Range("I132").Value = 2743.6858
$ws.Range("J132").Value = 4163.5
$ws.Range("K132").Value = 8231.057400000002
$ws.Range("L132").Value = 12490.5
$ws.Range("M132").Value = -5701.057400000002
$ws.Range("N132").Value = -17550.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 19833.334
$ws.Range("I13").Value = 27250
$ws.Range("J13").Value = 5000
$ws.Range("K13").Value = 27250
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = -27110
$ws.Range("N13").Value = -5280
$ws.Range("H40").Value = 3049.5405
$ws.Range("I40").Value = 1992.0476
$ws.Range("K40").Value = 1992.0476
$ws.Range("M40").Value = -1856.0476
$ws.Range("H46").Value = 1857
$ws.Range("I46").Value = 1299.6666
$ws.Range("J46").Value = 2275
$ws.Range("K46").Value = 1299.6666
$ws.Range("L46").Value = 2275
$ws.Range("M46").Value = -1111.6666
$ws.Range("N46").Value = -2651
$ws.Range("H61").Value = 6645.914
$ws.Range("I61").Value = 6453.4
$ws.Range("J61").Value = 7801
$ws.Range("K61").Value = 6453.4
$ws.Range("L61").Value = 7801
$ws.Range("M61").Value = -6251.4
$ws.Range("N61").Value = -8205
$ws.Range("H68").Value = 2364.8667
$ws.Range("I68").Value = 2226.7273
$ws.Range("J68").Value = 2744.75
$ws.Range("K68").Value = 2226.7273
$ws.Range("L68").Value = 2744.75
$ws.Range("M68").Value = -1477.7273
$ws.Range("N68").Value = -4242.75
$ws.Range("H71").Value = 2364.8667
$ws.Range("I71").Value = 2226.7273
$ws.Range("J71").Value = 2744.75
$ws.Range("K71").Value = 11133.6365
$ws.Range("L71").Value = 13723.75
$ws.Range("M71").Value = -7389.636500000001
$ws.Range("N71").Value = -21211.75
$ws.Range("H113").Value = 6645.914
$ws.Range("I113").Value = 6453.4
$ws.Range("J113").Value = 7801
$ws.Range("K113").Value = 6453.4
$ws.Range("L113").Value = 7801
$ws.Range("M113").Value = -4283.4
$ws.Range("N113").Value = -12141
$ws.Range("H132").Value = 3562.2942
$ws.Range("I132").Value = 2639.5454
$ws.Range("J132").Value = 5254
$ws.Range("K132").Value = 7918.6362
$ws.Range("L132").Value = 15762
$ws.Range("M132").Value = -5388.6362
$ws.Range("N132").Value = -20822

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 50000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H46").Value = 329933.34
$ws.Range("J46").Value = 329933.34
$ws.Range("L46").Value = 329933.34
$ws.Range("N46").Value = -330395.34
$ws.Range("H132").Value = 10002124
$ws.Range("I132").Value = 14707834
$ws.Range("J132").Value = 2488
$ws.Range("K132").Value = 44123502
$ws.Range("L132").Value = 7464
$ws.Range("M132").Value = -44120972
$ws.Range("N132").Value = -12524
$ws.Range("H134").Value = 329933.34
$ws.Range("J134").Value = 329933.34
$ws.Range("L134").Value = 989800.02
$ws.Range("N134").Value = -994870.02
$ws.Range("H136").Value = 6966253
$ws.Range("I136").Value = 9287993
$ws.Range("K136").Value = 27863979
$ws.Range("M136").Value = -27861429
